# Update the canonical terminology URLs and the generation date in the
# ValueSet-paysNaissance-vs workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value (row 8, column B) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Include #0 sheet: update the System URI (row 4, column B) ---
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R20-Pays/FHIR/TRE-R20-Pays"

# --- Include #1 sheet: update the System URI (row 4, column B) ---
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R268-PaysProvenanceISO/FHIR/TRE-R268-PaysProvenanceISO"

$wb.Save()
